$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.993.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.001.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.50%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.522"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.002.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.148"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.440"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000227"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.28%  "
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.488.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.953.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.999.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "458.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.679"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.13%  "
$ws.Range("E25").Value = "  -6.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.05%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.22%  "
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.108"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.07%  "
$ws.Range("E35").Value = "  -4.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0782"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.17%  "
$ws.Range("E41").Value = "  -12.03%  "
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.271"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.28%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0351"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.33%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "378.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -15.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.741.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "37.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.108"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.02%  "
